$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new row 6 to Sheet1, mirroring the existing rows' text-typed cells.
# A6 is an empty string (leading apostrophe forces Excel to keep it text
# instead of collapsing to a blank cell) and C6 is a numeric-looking string
# ("233") that also needs the text quote-prefix so it isn't stored as a
# number.
$ws.Range("A6").Value = "'"
$ws.Range("B6").Value = "احمد"
$ws.Range("C6").Value = "'233"
$ws.Range("D6").Value = "الصمود"
$ws.Range("E6").Value = "الرحلة 2"
$ws.Range("F6").Value = "C2"
$ws.Range("G6").Value = "IDRF"
$ws.Range("H6").Value = "٠٥‏/٠٥‏/٢٠٢٥ ٠٢:٢٧:١٧ م"
